# Apply the "Frais divers" hierarchical edits described in the commit:
#  - Entretien sheet: update the edited row's amount (B2) from a leftover
#    date placeholder to a real numeric amount, and bump B3.
#  - Maintenance sheet: rename the elevator service contract row for the
#    SourceIndex test.
#  - frais divers sheet (the combined/aggregated view, keyed by SourceIndex):
#    mirror the same edits so the hierarchy stays consistent.

$wb = $excel.ActiveWorkbook

# --- Entretien ---
$wsEntretien = $wb.Worksheets.Item("Entretien")
$wsEntretien.Range("B2").Value = 2500
$wsEntretien.Range("B3").Value = 750

# --- Maintenance ---
$wsMaintenance = $wb.Worksheets.Item("Maintenance")
$wsMaintenance.Range("A2").Value = "Modified by SourceIndex test"

# --- frais divers (aggregated sheet) ---
$wsFraisDivers = $wb.Worksheets.Item("frais divers")
$wsFraisDivers.Range("A2").Value = "Modified Description"
$wsFraisDivers.Range("B2").Value = 2500
$wsFraisDivers.Range("B3").Value = 750
$wsFraisDivers.Range("A8").Value = "Modified by SourceIndex test"
